$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 435, shifting rows 435:458 down to 436:459
$ws.Rows.Item(435).Insert()

# Fill the new row 435 with data (same metadata as the surrounding rows, new values per diff)
$ws.Range("A435").Value = 10
$ws.Range("B435").Value = "Vega Modelo de Temuco"
$ws.Range("C435").Value = "La Araucanía"
$ws.Range("D435").Value = 44585
$ws.Range("D435").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E435").Value = 9
$ws.Range("F435").Value = "Fruta"
$ws.Range("G435").Value = 100102
$ws.Range("H435").Value = "Cítricos"
$ws.Range("I435").Value = 100102004
$ws.Range("J435").Value = "Mandarina"
$ws.Range("K435").Value = "Murcott"
$ws.Range("L435").Value = "Primera"
$ws.Range("M435").Value = 150
$ws.Range("N435").Value = 14000
$ws.Range("O435").Value = 14000
$ws.Range("P435").Value = 14000
$ws.Range("Q435").Value = "$/bandeja 10 kilos"
$ws.Range("R435").Value = "Región de O'Higgins"
$ws.Range("S435").Value = 1400
$ws.Range("T435").Value = 10
